$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Pepino dulce" (Vega Modelo de
# Temuco). It belongs right before the existing row 178, so insert a blank
# row there; this pushes the old rows 178-280 down to 179-281 and extends
# the used range from A1:R280 to A1:R281.
$ws.Rows("178:178").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(178, 1).Value = 10
$ws.Cells.Item(178, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(178, 3).Value = "La Araucanía"
$ws.Cells.Item(178, 4).Value = 44813
$ws.Cells.Item(178, 5).Value = 9
$ws.Cells.Item(178, 6).Value = 100112043
$ws.Cells.Item(178, 7).Value = "Pepino dulce"
$ws.Cells.Item(178, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(178, 9).Value = "Primera"
$ws.Cells.Item(178, 10).Value = 100
$ws.Cells.Item(178, 11).Value = 18000
$ws.Cells.Item(178, 12).Value = 19000
$ws.Cells.Item(178, 13).Value = 18500
$ws.Cells.Item(178, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(178, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(178, 16).Value = 1028
$ws.Cells.Item(178, 17).Value = 18
$ws.Cells.Item(178, 18).Value = "Hortaliza"
